$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 80.109075200000007
$ws.Range("D3").Value = 80.109075200000007
$ws.Range("D4").Value = 80.109075200000007
$ws.Range("D5").Value = 80.109075200000007
